# Re-applies the upstream re-ordering of observation rows 2-5 and 12-14.
# The rows themselves were re-sequenced (same records, new row positions),
# so we snapshot every involved cell first and then write the snapshots
# back out to their new row positions. This avoids any intermediate
# overwrite since several of the moves form cycles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Superset of columns that are populated anywhere in the affected rows.
# (Y/AA - Startdatum/Slutdatum - are identical ("2026-01-22") across every
# row in both rotations, so they are intentionally left untouched: writing
# that same text back through Range.Value would let Excel's auto-detection
# reinterpret it as a date serial instead of leaving it as the original
# plain text.)
$cols = @("A","B","D","E","F","G","H","I","M","P","Q","R","S","T","U","V","W","Z","AB","AC","AD","AE","AG","AT","AW","AX","AY")

# Snapshot the current contents of every row we are about to touch.
$sourceRows = @(2,3,4,5,12,13,14)
$snapshot = @{}
foreach ($r in $sourceRows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# New row r gets the content that used to live in row map[r].
$rowMap = @{
    2  = 3
    3  = 5
    4  = 2
    5  = 4
    12 = 14
    13 = 12
    14 = 13
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $val = $rowData[$c]
        if ($val -eq $null) {
            $val = ""
        }
        $ws.Range("$c$destRow").Value = $val
    }
}
